$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "model_3_1_24"
$ws.Cells.Item(2, 2).Value = 0.3938518760383836
$ws.Cells.Item(2, 3).Value = 0.2950945914788321
$ws.Cells.Item(2, 4).Value = -0.3757217657599019
$ws.Cells.Item(2, 5).Value = 0.01489351509655346
$ws.Cells.Item(2, 6).Value = 0.6708273887634277
$ws.Cells.Item(2, 7).Value = 0.8525978922843933
$ws.Cells.Item(2, 8).Value = 1.38070011138916
$ws.Cells.Item(2, 9).Value = 1.101116538047791

$ws.Cells.Item(3, 1).Value = "model_3_1_23"
$ws.Cells.Item(3, 2).Value = 0.4022612072430602
$ws.Cells.Item(3, 3).Value = 0.2980495138885267
$ws.Cells.Item(3, 4).Value = -0.3125403898878389
$ws.Cells.Item(3, 5).Value = 0.04328221764052098
$ws.Cells.Item(3, 6).Value = 0.6615206599235535
$ws.Cells.Item(3, 7).Value = 0.8490238785743713
$ws.Cells.Item(3, 8).Value = 1.317290067672729
$ws.Cells.Item(3, 9).Value = 1.069384694099426

$ws.Cells.Item(4, 1).Value = "model_3_1_8"
$ws.Cells.Item(4, 2).Value = 0.4092409567790825
$ws.Cells.Item(4, 3).Value = 0.2451472509914322
$ws.Cells.Item(4, 4).Value = -0.09904503024820754
$ws.Cells.Item(4, 5).Value = 0.1031849232291139
$ws.Cells.Item(4, 6).Value = 0.6537961959838867
$ws.Cells.Item(4, 7).Value = 0.91301029920578
$ws.Cells.Item(4, 8).Value = 1.10302209854126
$ws.Cells.Item(4, 9).Value = 1.002427697181702

$ws.Cells.Item(5, 1).Value = "model_3_1_22"
$ws.Cells.Item(5, 2).Value = 0.4109394350476524
$ws.Cells.Item(5, 3).Value = 0.2985134748522186
$ws.Cells.Item(5, 4).Value = -0.2425344303859989
$ws.Cells.Item(5, 5).Value = 0.07312798547488653
$ws.Cells.Item(5, 6).Value = 0.65191650390625
$ws.Cells.Item(5, 7).Value = 0.8484627604484558
$ws.Cells.Item(5, 8).Value = 1.247030735015869
$ws.Cells.Item(5, 9).Value = 1.036024212837219

$ws.Cells.Item(6, 1).Value = "model_3_1_3"
$ws.Cells.Item(6, 2).Value = 0.4115407635989745
$ws.Cells.Item(6, 3).Value = 0.2117654556381081
$ws.Cells.Item(6, 4).Value = 0.08604093065382878
$ws.Cells.Item(6, 5).Value = 0.1622658909900693
$ws.Cells.Item(6, 6).Value = 0.651250958442688
$ws.Cells.Item(6, 7).Value = 0.9533863067626953
$ws.Cells.Item(6, 8).Value = 0.9172664880752563
$ws.Cells.Item(6, 9).Value = 0.9363890886306763

$ws.Cells.Item(7, 1).Value = "model_3_1_9"
$ws.Cells.Item(7, 2).Value = 0.4115733954997896
$ws.Cells.Item(7, 3).Value = 0.2558055473249926
$ws.Cells.Item(7, 4).Value = -0.1153395129717119
$ws.Cells.Item(7, 5).Value = 0.1024058506379429
$ws.Cells.Item(7, 6).Value = 0.6512148976325989
$ws.Cells.Item(7, 7).Value = 0.9001188278198242
$ws.Cells.Item(7, 8).Value = 1.119375586509705
$ws.Cells.Item(7, 9).Value = 1.00329852104187

$ws.Cells.Item(8, 1).Value = "model_3_1_2"
$ws.Cells.Item(8, 2).Value = 0.4127132118981602
$ws.Cells.Item(8, 3).Value = 0.2117347946403261
$ws.Cells.Item(8, 4).Value = 0.1009952487826695
$ws.Cells.Item(8, 5).Value = 0.1685672433128218
$ws.Cells.Item(8, 6).Value = 0.6499534249305725
$ws.Cells.Item(8, 7).Value = 0.9534233808517456
$ws.Cells.Item(8, 8).Value = 0.902258038520813
$ws.Cells.Item(8, 9).Value = 0.9293456077575684

$ws.Cells.Item(9, 1).Value = "model_3_1_21"
$ws.Cells.Item(9, 2).Value = 0.4166326723679018
$ws.Cells.Item(9, 3).Value = 0.2998827447789916
$ws.Cells.Item(9, 4).Value = -0.200324558497482
$ws.Cells.Item(9, 5).Value = 0.09174742806755087
$ws.Cells.Item(9, 6).Value = 0.6456156969070435
$ws.Cells.Item(9, 7).Value = 0.8468065857887268
$ws.Cells.Item(9, 8).Value = 1.204668283462524
$ws.Cells.Item(9, 9).Value = 1.015212059020996

$ws.Cells.Item(10, 1).Value = "model_3_1_13"
$ws.Cells.Item(10, 2).Value = 0.4191467324253911
$ws.Cells.Item(10, 3).Value = 0.2863176500543526
$ws.Cells.Item(10, 4).Value = -0.1408249454200614
$ws.Cells.Item(10, 5).Value = 0.1091168931855169
$ws.Cells.Item(10, 6).Value = 0.642833411693573
$ws.Cells.Item(10, 7).Value = 0.863213837146759
$ws.Cells.Item(10, 8).Value = 1.14495325088501
$ws.Cells.Item(10, 9).Value = 0.9957970976829529

$ws.Cells.Item(11, 1).Value = "model_3_1_12"
$ws.Cells.Item(11, 2).Value = 0.4195577879652058
$ws.Cells.Item(11, 3).Value = 0.287364119712481
$ws.Cells.Item(11, 4).Value = -0.1355751026217755
$ws.Cells.Item(11, 5).Value = 0.1119346122398783
$ws.Cells.Item(11, 6).Value = 0.6423785090446472
$ws.Cells.Item(11, 7).Value = 0.8619481325149536
$ws.Cells.Item(11, 8).Value = 1.139684438705444
$ws.Cells.Item(11, 9).Value = 0.9926475286483765

$ws.Cells.Item(12, 1).Value = "model_3_1_20"
$ws.Cells.Item(12, 2).Value = 0.4198423089939615
$ws.Cells.Item(12, 3).Value = 0.298393953625926
$ws.Cells.Item(12, 4).Value = -0.1716195180926749
$ws.Cells.Item(12, 5).Value = 0.1030235618510323
$ws.Cells.Item(12, 6).Value = 0.6420636177062988
$ws.Cells.Item(12, 7).Value = 0.8486072421073914
$ws.Cells.Item(12, 8).Value = 1.175859212875366
$ws.Cells.Item(12, 9).Value = 1.002608060836792

$ws.Cells.Item(13, 1).Value = "model_3_1_19"
$ws.Cells.Item(13, 2).Value = 0.4205764848778644
$ws.Cells.Item(13, 3).Value = 0.2983616593008557
$ws.Cells.Item(13, 4).Value = -0.1649739238123886
$ws.Cells.Item(13, 5).Value = 0.1058129769575293
$ws.Cells.Item(13, 6).Value = 0.6412511467933655
$ws.Cells.Item(13, 7).Value = 0.8486464023590088
$ws.Cells.Item(13, 8).Value = 1.169189691543579
$ws.Cells.Item(13, 9).Value = 0.9994900822639465

$ws.Cells.Item(14, 1).Value = "model_3_1_11"
$ws.Cells.Item(14, 2).Value = 0.4218613314113581
$ws.Cells.Item(14, 3).Value = 0.2806334138040402
$ws.Cells.Item(14, 4).Value = -0.1047676341566859
$ws.Cells.Item(14, 5).Value = 0.1210957950180442
$ws.Cells.Item(14, 6).Value = 0.639829158782959
$ws.Cells.Item(14, 7).Value = 0.8700889945030212
$ws.Cells.Item(14, 8).Value = 1.108765482902527
$ws.Cells.Item(14, 9).Value = 0.9824075102806091

$ws.Cells.Item(15, 1).Value = "model_3_1_10"
$ws.Cells.Item(15, 2).Value = 0.4218628268982997
$ws.Cells.Item(15, 3).Value = 0.2816693669345343
$ws.Cells.Item(15, 4).Value = -0.1050496007442152
$ws.Cells.Item(15, 5).Value = 0.1215695019472159
$ws.Cells.Item(15, 6).Value = 0.6398274302482605
$ws.Cells.Item(15, 7).Value = 0.8688360452651978
$ws.Cells.Item(15, 8).Value = 1.10904848575592
$ws.Cells.Item(15, 9).Value = 0.9818779230117798

$ws.Cells.Item(16, 1).Value = "model_3_1_15"
$ws.Cells.Item(16, 2).Value = 0.4220057600055157
$ws.Cells.Item(16, 3).Value = 0.2858328517151734
$ws.Cells.Item(16, 4).Value = -0.12234379034855
$ws.Cells.Item(16, 5).Value = 0.1166482454657525
$ws.Cells.Item(16, 6).Value = 0.6396693587303162
$ws.Cells.Item(16, 7).Value = 0.8638002872467041
$ws.Cells.Item(16, 8).Value = 1.126405239105225
$ws.Cells.Item(16, 9).Value = 0.9873788356781006

$ws.Cells.Item(17, 1).Value = "model_3_1_14"
$ws.Cells.Item(17, 2).Value = 0.4228577919836178
$ws.Cells.Item(17, 3).Value = 0.2851611435595762
$ws.Cells.Item(17, 4).Value = -0.1138181600271284
$ws.Cells.Item(17, 5).Value = 0.1198658463882191
$ws.Cells.Item(17, 6).Value = 0.6387263536453247
$ws.Cells.Item(17, 7).Value = 0.8646126985549927
$ws.Cells.Item(17, 8).Value = 1.117848753929138
$ws.Cells.Item(17, 9).Value = 0.9837823510169983

$ws.Cells.Item(18, 1).Value = "model_3_1_18"
$ws.Cells.Item(18, 2).Value = 0.4229449346633816
$ws.Cells.Item(18, 3).Value = 0.2966990705263038
$ws.Cells.Item(18, 4).Value = -0.1411077750290963
$ws.Cells.Item(18, 5).Value = 0.1149447470356773
$ws.Cells.Item(18, 6).Value = 0.6386299133300781
$ws.Cells.Item(18, 7).Value = 0.8506572842597961
$ws.Cells.Item(18, 8).Value = 1.14523708820343
$ws.Cells.Item(18, 9).Value = 0.9892830848693848

$ws.Cells.Item(19, 1).Value = "model_3_1_4"
$ws.Cells.Item(19, 2).Value = 0.4230840833677169
$ws.Cells.Item(19, 3).Value = 0.2208521905157045
$ws.Cells.Item(19, 4).Value = 0.1088730531513687
$ws.Cells.Item(19, 5).Value = 0.1771190752729923
$ws.Cells.Item(19, 6).Value = 0.6384759545326233
$ws.Cells.Item(19, 7).Value = 0.9423957467079163
$ws.Cells.Item(19, 8).Value = 0.8943517208099365
$ws.Cells.Item(19, 9).Value = 0.9197867512702942

$ws.Cells.Item(20, 1).Value = "model_3_1_5"
$ws.Cells.Item(20, 2).Value = 0.4239896329582331
$ws.Cells.Item(20, 3).Value = 0.22157972355688
$ws.Cells.Item(20, 4).Value = 0.1107076003281058
$ws.Cells.Item(20, 5).Value = 0.1783107561288255
$ws.Cells.Item(20, 6).Value = 0.6374737024307251
$ws.Cells.Item(20, 7).Value = 0.9415156841278076
$ws.Cells.Item(20, 8).Value = 0.8925105333328247
$ws.Cells.Item(20, 9).Value = 0.918454647064209

$ws.Cells.Item(21, 1).Value = "model_3_1_1"
$ws.Cells.Item(21, 2).Value = 0.4241479393855609
$ws.Cells.Item(21, 3).Value = 0.2378741848211962
$ws.Cells.Item(21, 4).Value = 0.1608706814201163
$ws.Cells.Item(21, 5).Value = 0.2088410401937679
$ws.Cells.Item(21, 6).Value = 0.6372985243797302
$ws.Cells.Item(21, 7).Value = 0.9218072891235352
$ws.Cells.Item(21, 8).Value = 0.8421658873558044
$ws.Cells.Item(21, 9).Value = 0.884329080581665

$ws.Cells.Item(22, 1).Value = "model_3_1_17"
$ws.Cells.Item(22, 2).Value = 0.4241853367472079
$ws.Cells.Item(22, 3).Value = 0.2976635572887913
$ws.Cells.Item(22, 4).Value = -0.1308888562054242
$ws.Cells.Item(22, 5).Value = 0.1198149387725937
$ws.Cells.Item(22, 6).Value = 0.6372570991516113
$ws.Cells.Item(22, 7).Value = 0.8494907021522522
$ws.Cells.Item(22, 8).Value = 1.134981155395508
$ws.Cells.Item(22, 9).Value = 0.9838391542434692

$ws.Cells.Item(23, 1).Value = "model_3_1_0"
$ws.Cells.Item(23, 2).Value = 0.4244802528636702
$ws.Cells.Item(23, 3).Value = 0.2606389186613304
$ws.Cells.Item(23, 4).Value = 0.1924712546504692
$ws.Cells.Item(23, 5).Value = 0.2352347329204794
$ws.Cells.Item(23, 6).Value = 0.6369307637214661
$ws.Cells.Item(23, 7).Value = 0.8942728042602539
$ws.Cells.Item(23, 8).Value = 0.8104510307312012
$ws.Cells.Item(23, 9).Value = 0.8548271059989929

$ws.Cells.Item(24, 1).Value = "model_3_1_16"
$ws.Cells.Item(24, 2).Value = 0.4270213038272854
$ws.Cells.Item(24, 3).Value = 0.2946662127551453
$ws.Cells.Item(24, 4).Value = -0.09954972837390375
$ws.Cells.Item(24, 5).Value = 0.1313397561789065
$ws.Cells.Item(24, 6).Value = 0.6341186165809631
$ws.Cells.Item(24, 7).Value = 0.8531160354614258
$ws.Cells.Item(24, 8).Value = 1.103528738021851
$ws.Cells.Item(24, 9).Value = 0.9709572792053223

$ws.Cells.Item(25, 1).Value = "model_3_1_6"
$ws.Cells.Item(25, 2).Value = 0.434378637612327
$ws.Cells.Item(25, 3).Value = 0.2329395146073177
$ws.Cells.Item(25, 4).Value = 0.1538724210486149
$ws.Cells.Item(25, 5).Value = 0.2030573382656511
$ws.Cells.Item(25, 6).Value = 0.6259761452674866
$ws.Cells.Item(25, 7).Value = 0.9277758002281189
$ws.Cells.Item(25, 8).Value = 0.8491895198822021
$ws.Cells.Item(25, 9).Value = 0.8907938599586487

$ws.Cells.Item(26, 1).Value = "model_3_1_7"
$ws.Cells.Item(26, 2).Value = 0.44020897982787
$ws.Cells.Item(26, 3).Value = 0.2403448007179627
$ws.Cells.Item(26, 4).Value = 0.1618025437671115
$ws.Cells.Item(26, 5).Value = 0.2106501843699528
$ws.Cells.Item(26, 6).Value = 0.6195237040519714
$ws.Cells.Item(26, 7).Value = 0.9188190102577209
$ws.Cells.Item(26, 8).Value = 0.8412306308746338
$ws.Cells.Item(26, 9).Value = 0.8823068141937256

